$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D hold price text that can look like numbers
# (e.g. "408.45"), so force the cell to Text format before assigning the
# string value, then restore the default "Normal" style so the saved file
# keeps the same (unstyled) look as the rest of the sheet.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '62.129.31'
$ws.Range('E2').Value = '  +1.83%  '
Set-TextValue $ws.Range('D3') '3.435.49'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '408.45'
$ws.Range('E5').Value = '  +0.78%  '
Set-TextValue $ws.Range('D6') '128.66'
$ws.Range('E6').Value = '  -3.26%  '
Set-TextValue $ws.Range('D7') '0.631'
$ws.Range('E7').Value = '  +7.17%  '
$ws.Range('E8').Value = '  -0.15%  '
Set-TextValue $ws.Range('D9') '0.749'
$ws.Range('E9').Value = '  +11.87%  '
$ws.Range('E10').Value = '  +17.86%  '
Set-TextValue $ws.Range('D11') '43.42'
$ws.Range('E11').Value = '  +3.04%  '
Set-TextValue $ws.Range('D12') '0.140'
$ws.Range('E12').Value = '  -0.33%  '
Set-TextValue $ws.Range('D13') '8.73'
$ws.Range('E13').Value = '  +5.15%  '
Set-TextValue $ws.Range('D14') '20.27'
$ws.Range('E14').Value = '  +3.70%  '
$ws.Range('E15').Value = '  +56.02%  '
Set-TextValue $ws.Range('D16') '3.446.78'
$ws.Range('E16').Value = '  +2.81%  '
Set-TextValue $ws.Range('D17') '62.097.32'
$ws.Range('E18').Value = '  +2.86%  '
Set-TextValue $ws.Range('D19') '11.46'
$ws.Range('E19').Value = '  +3.25%  '
Set-TextValue $ws.Range('D20') '371.49'
$ws.Range('E20').Value = '  +22.20%  '
Set-TextValue $ws.Range('D21') '87.20'
$ws.Range('E21').Value = '  +4.31%  '
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('E23').Value = '  +4.44%  '
$ws.Range('E24').Value = '  +2.22%  '
Set-TextValue $ws.Range('D25') '31.69'
$ws.Range('E25').Value = '  +8.03%  '
$ws.Range('E26').Value = '  +0.60%  '
Set-TextValue $ws.Range('D27') '8.46'
$ws.Range('E27').Value = '  +2.00%  '
Set-TextValue $ws.Range('D28') '7.69'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('E29').Value = '  +10.45%  '
Set-TextValue $ws.Range('D30') '43.96'
$ws.Range('E30').Value = '  +6.81%  '
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('E32').Value = '  +0.20%  '
Set-TextValue $ws.Range('D33') '11.80'
$ws.Range('E33').Value = '  +4.62%  '
$ws.Range('E34').Value = '  -0.07%  '
Set-TextValue $ws.Range('D35') '0.0492'
$ws.Range('E35').Value = '  +2.95%  '
Set-TextValue $ws.Range('D36') '52.10'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D40') '0.132'
$ws.Range('E40').Value = '  +7.03%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D41') '0.313'
$ws.Range('E41').Value = '  +9.47%  '
Set-TextValue $ws.Range('D42') '143.15'
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  +0.80%  '
Set-TextValue $ws.Range('D45') '16.64'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +4.53%  '
Set-TextValue $ws.Range('D47') '21.72'
$ws.Range('E47').Value = '  +1.35%  '
Set-TextValue $ws.Range('D48') '2.112.93'
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('E51').Value = '  +6.92%  '
